# Add accent6 theme-color fill to the runs of reviewer point "5." (both the
# English and the Korean text boxes) on slide 2, explaining KPF-BERT,
# ETRI-ELECTRA and ETRI-RoBERTa.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(2)

# --- English content placeholder (shape 2) -------------------------------
$shpEn = $s.Shapes.Item(2)
$trEn = $shpEn.TextFrame.TextRange
$fullEn = $trEn.Text
$targetEn = "5. The Author should provide some more information on three pre-trained language models, KPF-BERT, ETRI-ELECTRA, and ETRI-RoBERTa, which they are using to fine-tune their re-ranking model."
$idxEn = $fullEn.IndexOf($targetEn)
if ($idxEn -ge 0) {
    $subEn = $trEn.Characters($idxEn + 1, $targetEn.Length)
    $subEn.Font.Color.ObjectThemeColor = 10  # msoThemeColorAccent6
}

# --- Korean content placeholder (shape 3) --------------------------------
$shpKo = $s.Shapes.Item(3)
$trKo = $shpKo.TextFrame.TextRange
$fullKo = $trKo.Text
$targetKo = "5. 저자는 자신의 재순위화 모델을 세밀하게 조정하기 위해 사용하고 있는 사전 학습된 언어 모델 세 가지(KPF-BERT, ETRI-ELECTRA, ETRI-RoBERTa)에 대한 추가 정보를 제공해야 합니다."
$idxKo = $fullKo.IndexOf($targetKo)
if ($idxKo -ge 0) {
    $subKo = $trKo.Characters($idxKo + 1, $targetKo.Length)
    $subKo.Font.Color.ObjectThemeColor = 10  # msoThemeColorAccent6
}
